# Add 2022-Q3 data:
#  1. Insert a new summary row (row 2) in "总计" with the 2022-Q3 totals,
#     shifting the existing quarters down by one row.
#  2. Insert a new worksheet "2022-Q3" (cloned from "2022-Q2" so it keeps
#     identical formatting) right after "总计" and fill it with the new
#     per-fund holding data for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" overview sheet: insert new row 2 for 2022-Q3
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows(2).Insert()
$totals.Range("B2:D2").ClearFormats()

# Copy the numbering-column style (bold/centered/bordered) from the row
# below (the shifted former row 2) so the new index cell matches.
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 4
$totals.Range("D2").Value = 1.09

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet with fund-holding detail rows
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($wb.Worksheets.Item("总计").Next())
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Row 2 (fund 001643) - code/name unchanged, figures updated
$q3.Range("D2:G2").NumberFormat = "@"
$q3.Range("D2").Value = "15.68"
$q3.Range("E2").Value = "94.47"
$q3.Range("F2").Value = "3.70"
$q3.Range("G2").Value = "0.5802"
$q3.Range("H2").Value = 10

# Row 3 (fund 001644) - code/name unchanged, figures updated
$q3.Range("D3:G3").NumberFormat = "@"
$q3.Range("D3").Value = "8.27"
$q3.Range("E3").Value = "94.47"
$q3.Range("F3").Value = "3.70"
$q3.Range("G3").Value = "0.3060"
$q3.Range("H3").Value = 10

# New row 4 (fund 014575)
$q3.Range("A4").Value = 2
$q3.Range("B4:G4").NumberFormat = "@"
$q3.Range("B4").Value = "014575"
$q3.Range("C4").Value = "鑫元清洁能源混合C"
$q3.Range("D4").Value = "1.77"
$q3.Range("E4").Value = "93.05"
$q3.Range("F4").Value = "7.99"
$q3.Range("G4").Value = "0.1414"
$q3.Range("H4").Value = 7

# New row 5 (fund 014574)
$q3.Range("A5").Value = 3
$q3.Range("B5:G5").NumberFormat = "@"
$q3.Range("B5").Value = "014574"
$q3.Range("C5").Value = "鑫元清洁能源混合A"
$q3.Range("D5").Value = "0.82"
$q3.Range("E5").Value = "93.05"
$q3.Range("F5").Value = "7.99"
$q3.Range("G5").Value = "0.0655"
$q3.Range("H5").Value = 7

# Match the bold/centered/bordered numbering-column style used by rows 2-3
$q3.Range("A2").Copy()
$q3.Range("A4:A5").PasteSpecial(-4122)
$q3.Range("A4").Value = 2
$q3.Range("A5").Value = 3
